{"js": "// Replace \"monitored\" with \"inspected\" in the sentence about packet\n// captures: \"Traffic will only be monitored or stored (e.g. packet\n// captured)\" -> \"Traffic will only be inspected or stored (e.g. packet\n// captured)\".\nconst searchResults = context.document.body.search(\"monitored or stored (e.g. packet captured)\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items/text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find target text to update.\");\n}\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  const found = searchResults.items[i];\n  found.insertText(\"inspected or stored (e.g. packet captured)\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Fix the language in the last change: s/monitored/inspected/ for the\n# sentence about packet captures.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$found = $find.Execute(\n    \"monitored\",            # FindText\n    $true,                  # MatchCase\n    $false,                 # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    \"wdFindContinue\",       # Wrap\n    $false,                 # Format\n    \"inspected\",            # ReplaceWith\n    \"wdReplaceOne\"          # Replace\n)\n\nif (-not $found) {\n    throw \"Could not find target text 'monitored' to update.\"\n}\n"}
